$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-10 from 45183 (2023-09-14)
# to 45184 (2023-09-15), matching the diff.
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45184
}
